# [Kadastro App] Yeni kayit eklendi: 3008
#
# Appends a new record (row 66) to both the master "Kayitlar" log sheet and
# the "Erdemli" unit sheet (which mirrors Kayitlar's Erdemli-filtered rows).
# All columns in this workbook are stored as text, so every cell on the new
# row is forced to a text number-format before the value is written - this
# keeps numeric-looking values ("3008", "1") and the date ("2025-09-11") as
# literal text instead of being auto-converted to numbers/dates.

$wb = $excel.ActiveWorkbook

$newRow = @{
    A = "3008"
    B = "2025-09-11"
    C = "Erdemli"
    D = "1"
    E = "CİNS DEĞ."
    F = "CEMAL TİMUROĞLU (K.Teknisyeni), ÖZKAN AKBAŞ (Mühendis)"
}

$sheetNames = @("Kayitlar", "Erdemli")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $targetRow = $ws.UsedRange.Rows.Count + 1

    $rowRange = $ws.Range("A" + $targetRow + ":F" + $targetRow)
    $rowRange.NumberFormat = "@"

    $ws.Cells.Item($targetRow, 1).Value = $newRow.A
    $ws.Cells.Item($targetRow, 2).Value = $newRow.B
    $ws.Cells.Item($targetRow, 3).Value = $newRow.C
    $ws.Cells.Item($targetRow, 4).Value = $newRow.D
    $ws.Cells.Item($targetRow, 5).Value = $newRow.E
    $ws.Cells.Item($targetRow, 6).Value = $newRow.F

    # Mirrors Excel's "Ignore Error" action for numbers stored as text over
    # the whole used range (matches the pre-existing ignoredError on this
    # sheet, just widened to cover the newly appended row).
    $usedRange = $ws.Range("A1:F" + $targetRow)
    $usedRange.Errors.Item(1).Ignore = $true
}
